$wb = $excel.ActiveWorkbook

$wsMean = $wb.Worksheets.Item("Stations_Mean")
$wsStd = $wb.Worksheets.Item("Stations_Std")
$wsCV = $wb.Worksheets.Item("Stations_CV")

# Stations_Mean
$wsMean.Range("B3").Value = 4.179417073170732
$wsMean.Range("C3").Value = 27.44548524590163
$wsMean.Range("B4").Value = 4.818949350649348
$wsMean.Range("C4").Value = 26.81466027397261
$wsMean.Range("B5").Value = 4.161407594936702
$wsMean.Range("C5").Value = 26.73234054054054
$wsMean.Range("B6").Value = 4.24243783783784
$wsMean.Range("C6").Value = 26.74118630136986
$wsMean.Range("B7").Value = 4.426701449275361
$wsMean.Range("C7").Value = 26.16007397260275
$wsMean.Range("B8").Value = 4.471012328767126
$wsMean.Range("C8").Value = 25.90625479452053
$wsMean.Range("B9").Value = 4.794501315789475
$wsMean.Range("C9").Value = 25.47709999999999
$wsMean.Range("B10").Value = 6.017502702702703
$wsMean.Range("C10").Value = 25.8668366197183
$wsMean.Range("B11").Value = 5.187090909090905
$wsMean.Range("C11").Value = 30.97843116883117
$wsMean.Range("B12").Value = 4.643550000000002
$wsMean.Range("C12").Value = 33.74741408450704
$wsMean.Range("B13").Value = 4.714093243243243
$wsMean.Range("C13").Value = 36.10643661971831
$wsMean.Range("B14").Value = 5.829694520547945
$wsMean.Range("C14").Value = 38.14268169014085

# Stations_Std
$wsStd.Range("B3").Value = 0.220873811131656
$wsStd.Range("C3").Value = 0.8380774200091998
$wsStd.Range("B4").Value = 0.3055942475103465
$wsStd.Range("C4").Value = 0.7171445489539521
$wsStd.Range("B5").Value = 0.0556212899013978
$wsStd.Range("C5").Value = 0.1953421465643332
$wsStd.Range("B6").Value = 0.1631128751985414
$wsStd.Range("C6").Value = 0.2782024968510989
$wsStd.Range("B7").Value = 0.08913679753244069
$wsStd.Range("C7").Value = 0.3165358352661047
$wsStd.Range("B8").Value = 0.063355844563346
$wsStd.Range("C8").Value = 0.2299552209073035
$wsStd.Range("B9").Value = 0.186891195079505
$wsStd.Range("C9").Value = 0.3680670150019801
$wsStd.Range("B10").Value = 0.7941893168203027
$wsStd.Range("C10").Value = 0.4861261427241282
$wsStd.Range("B11").Value = 0.128215334067155
$wsStd.Range("C11").Value = 0.6854042626262595
$wsStd.Range("B12").Value = 0.05141226591929673
$wsStd.Range("C12").Value = 0.5288059523373055
$wsStd.Range("B13").Value = 0.07928776846206127
$wsStd.Range("C13").Value = 0.5422719110344759
$wsStd.Range("B14").Value = 0.1343797850728497
$wsStd.Range("C14").Value = 0.5369750561671127

# Stations_CV
$wsCV.Range("B3").Value = 5.284799465206023
$wsCV.Range("C3").Value = 3.053607587915932
$wsCV.Range("B4").Value = 6.341511920418255
$wsCV.Range("C4").Value = 2.674449504959947
$wsCV.Range("B5").Value = 1.336597981151227
$wsCV.Range("C5").Value = 0.730733421071342
$wsCV.Range("B6").Value = 3.844791165677326
$wsCV.Range("C6").Value = 1.04035211346195
$wsCV.Range("B7").Value = 2.013616652350299
$wsCV.Range("C7").Value = 1.209995948779084
$wsCV.Range("B8").Value = 1.417035783053103
$wsCV.Range("C8").Value = 0.8876436317454179
$wsCV.Range("B9").Value = 3.898031990606118
$wsCV.Range("C9").Value = 1.444697453799609
$wsCV.Range("B10").Value = 13.19798853540357
$wsCV.Range("C10").Value = 1.879341296622077
$wsCV.Range("B11").Value = 2.471815827296271
$wsCV.Range("C11").Value = 2.212520895234606
$wsCV.Range("B12").Value = 1.107175887398579
$wsCV.Range("C12").Value = 1.566952510829779
$wsCV.Range("B13").Value = 1.681930423750214
$wsCV.Range("C13").Value = 1.501870474635352
$wsCV.Range("B14").Value = 2.305091366266291
$wsCV.Range("C14").Value = 1.407806248468132
